$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 update: PT Prima Tunggal Mandiri (Shell)
$ws.Range("B8").Value = 44843
$ws.Range("C8").Value = 45939
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 160000000
$ws.Range("G8").Value = 160000000
$ws.Range("H8").Value = "Split Per Year"
